# Generate Report for Handoff
# Inserts a new tracked file "b394b7d0-767e-40fa-887b-60022579535f" as a new
# row between the existing "1870b237-..." and "e586e662-..." rows on all
# three worksheets (Overview, zh-cn, de-de), mirroring the existing layout,
# hyperlinks, and styling already used for those two rows.

$wb = $excel.ActiveWorkbook

$mdBase    = "https://github.com/OpenLocalizationTest/oltest/blob/31cf5cc15b0537210cec336d357bad9cb6dfb039/e2e"
$zhXlfBase = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/93739e22065b8d10eb04cbb29b25b8d6955d267e/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht"
$deXlfBase = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/cd2ac63657323f70b5777c70048ddc189f749d5f/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht"

$newFile = "b394b7d0-767e-40fa-887b-60022579535f"
$newSha  = "677eb8697dcf2f5b6af4018dddeb84f55a95f841"
$newMdUrl = "$mdBase/$newFile.md"

$oldMdUrl1 = "$mdBase/1870b237-848c-411b-b5ce-ae973b78bbb4.md"
$oldMdUrl2 = "$mdBase/e586e662-88af-4177-8c06-d4c87cbe033a.md"

# ---------------------------------------------------------------------------
# Sheet 1: Overview
# ---------------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("Overview")

$ws1.Rows("2").Copy()
$ws1.Rows("3").Insert()

$ws1.Range("A3").Value2 = "$newFile.md"
$ws1.Range("B3").Value2 = "Ready for handoff"
$ws1.Range("C3").Value2 = "Ready for handoff"
$ws1.Range("D3").Value2 = "2016-47-13 12:47:37"

$ws1.Range("A1").Hyperlinks.Delete()
$ws1.Hyperlinks.Add($ws1.Range("A2"), $oldMdUrl1, "", "", "1870b237-848c-411b-b5ce-ae973b78bbb4.md")
$ws1.Hyperlinks.Add($ws1.Range("A3"), $newMdUrl, "", "", "$newFile.md")
$ws1.Hyperlinks.Add($ws1.Range("A4"), $oldMdUrl2, "", "", "e586e662-88af-4177-8c06-d4c87cbe033a.md")

# ---------------------------------------------------------------------------
# Sheet 2: zh-cn
# ---------------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("zh-cn")

$ws2.Rows("2").Copy()
$ws2.Rows("3").Insert()

$ws2.Range("A3").Value2 = "$newFile.md"
$ws2.Range("B3").Value2 = ".md"
$ws2.Range("C3").Value2 = "Ready for handoff"
$ws2.Range("D3").Value2 = "$newFile.$newSha.zh-cn.xlf"
$ws2.Range("E3").Value2 = "2016-03-13 12:47:34"
$ws2.Range("H3").Value2 = "0001-01-01 00:00:00"
$ws2.Range("I3").Value2 = "Include"

$ws2.Range("A1").Hyperlinks.Delete()
$ws2.Hyperlinks.Add($ws2.Range("A2"), $oldMdUrl1, "", "", "1870b237-848c-411b-b5ce-ae973b78bbb4.md")
$ws2.Hyperlinks.Add($ws2.Range("B2"), $oldMdUrl1, "", "", ".md")
$ws2.Hyperlinks.Add($ws2.Range("D2"), "$zhXlfBase/1870b237-848c-411b-b5ce-ae973b78bbb4.9b70cfcef9f82634c21c0abd500437bbeb55aa4a.zh-cn.xlf", "", "", "1870b237-848c-411b-b5ce-ae973b78bbb4.9b70cfcef9f82634c21c0abd500437bbeb55aa4a.zh-cn.xlf")

$ws2.Hyperlinks.Add($ws2.Range("A3"), $newMdUrl, "", "", "$newFile.md")
$ws2.Hyperlinks.Add($ws2.Range("B3"), $newMdUrl, "", "", ".md")
$ws2.Hyperlinks.Add($ws2.Range("D3"), "$zhXlfBase/$newFile.$newSha.zh-cn.xlf", "", "", "$newFile.$newSha.zh-cn.xlf")

$ws2.Hyperlinks.Add($ws2.Range("A4"), $oldMdUrl2, "", "", "e586e662-88af-4177-8c06-d4c87cbe033a.md")
$ws2.Hyperlinks.Add($ws2.Range("B4"), $oldMdUrl2, "", "", ".md")
$ws2.Hyperlinks.Add($ws2.Range("D4"), "$zhXlfBase/e586e662-88af-4177-8c06-d4c87cbe033a.ec3ec9aa71759309e2f3c81a417dd408c97b0e10.zh-cn.xlf", "", "", "e586e662-88af-4177-8c06-d4c87cbe033a.ec3ec9aa71759309e2f3c81a417dd408c97b0e10.zh-cn.xlf")

# ---------------------------------------------------------------------------
# Sheet 3: de-de
# ---------------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("de-de")

$ws3.Rows("2").Copy()
$ws3.Rows("3").Insert()

$ws3.Range("A3").Value2 = "$newFile.md"
$ws3.Range("B3").Value2 = ".md"
$ws3.Range("C3").Value2 = "Ready for handoff"
$ws3.Range("D3").Value2 = "$newFile.$newSha.de-de.xlf"
$ws3.Range("E3").Value2 = "2016-03-13 12:47:37"
$ws3.Range("H3").Value2 = "0001-01-01 00:00:00"
$ws3.Range("I3").Value2 = "Include"

$ws3.Range("A1").Hyperlinks.Delete()
$ws3.Hyperlinks.Add($ws3.Range("A2"), $oldMdUrl1, "", "", "1870b237-848c-411b-b5ce-ae973b78bbb4.md")
$ws3.Hyperlinks.Add($ws3.Range("B2"), $oldMdUrl1, "", "", ".md")
$ws3.Hyperlinks.Add($ws3.Range("D2"), "$deXlfBase/1870b237-848c-411b-b5ce-ae973b78bbb4.9b70cfcef9f82634c21c0abd500437bbeb55aa4a.de-de.xlf", "", "", "1870b237-848c-411b-b5ce-ae973b78bbb4.9b70cfcef9f82634c21c0abd500437bbeb55aa4a.de-de.xlf")

$ws3.Hyperlinks.Add($ws3.Range("A3"), $newMdUrl, "", "", "$newFile.md")
$ws3.Hyperlinks.Add($ws3.Range("B3"), $newMdUrl, "", "", ".md")
$ws3.Hyperlinks.Add($ws3.Range("D3"), "$deXlfBase/$newFile.$newSha.de-de.xlf", "", "", "$newFile.$newSha.de-de.xlf")

$ws3.Hyperlinks.Add($ws3.Range("A4"), $oldMdUrl2, "", "", "e586e662-88af-4177-8c06-d4c87cbe033a.md")
$ws3.Hyperlinks.Add($ws3.Range("B4"), $oldMdUrl2, "", "", ".md")
$ws3.Hyperlinks.Add($ws3.Range("D4"), "$deXlfBase/e586e662-88af-4177-8c06-d4c87cbe033a.ec3ec9aa71759309e2f3c81a417dd408c97b0e10.de-de.xlf", "", "", "e586e662-88af-4177-8c06-d4c87cbe033a.ec3ec9aa71759309e2f3c81a417dd408c97b0e10.de-de.xlf")

Write-Host "Report generated for handoff: inserted $newFile rows on Overview, zh-cn, de-de."
